$wb = $excel.ActiveWorkbook

# Rename "sample" -> "RestAssured"
$ws = $wb.Worksheets.Item("sample")
$ws.Name = "RestAssured"

# Populate header row + sample data row
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "isbn "
$ws.Range("D1").Value = "aisle"
$ws.Range("E1").Value = "author"

$ws.Range("A2").Value = "Rest AddBook "
$ws.Range("B2").Value = "abcd "
$ws.Range("C2").Value = "abcd"
$ws.Range("D2").Value = 5765
$ws.Range("E2").Value = "Barman"

# A1 / C2 pick up the "no explicit color" Arial font (matching the testdata
# sheet's original unstyled font) with no border - grab it by copying the
# existing format from the testdata sheet and stripping the border back off.
$ws2 = $wb.Worksheets.Item("testdata")
$ws2.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Borders.LineStyle = -4142

$ws2.Range("A1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Borders.LineStyle = -4142

# Remaining new cells pick up the Arial + automatic(theme1) text colour font.
$ws.Range("B1:E1").Font.ThemeColor = 1
$ws.Range("A2:B2").Font.ThemeColor = 1
$ws.Range("D2:E2").Font.ThemeColor = 1

# The existing "testdata" sheet's Name/Value columns move to the same new
# Arial + theme1 font (keeping their border).
$ws2.Range("B1:C6").Font.ThemeColor = 1
